$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column B's width so the newly inserted column C can match it.
$bWidth = $ws.Columns("B").ColumnWidth

# Insert a new column before column C; this shifts C:X to D:Y for every row.
$ws.Columns("C").Insert()

# The insert carries B's style into the new C cells for every row that had
# a left-neighbour cell - but rows 6 and 10 should end up with no cell at
# all in the new column, so drop what Insert() auto-created there.
$ws.Range("C6").Clear()
$ws.Range("C10").Clear()

# Give the new column the same (approximate) width as column B.
$ws.Columns("C").ColumnWidth = $bWidth

# Header text for the new column, with its own distinct formatting
# (Arial 11, dark grey) rather than the bold header style copied from B1.
$ws.Range("C1").ClearFormats()
$ws.Range("C1").Value = "Project Number 2"
$ws.Range("C1").Font.Name = "Arial"
$ws.Range("C1").Font.Color = 3355443

$ws.Range("C1").Select()
